# Remove the duplicate/erroneous row 128 ("キウイフルーツ" post), which shifts
# all subsequent rows (129..289) up by one (128..288).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Delete()
